$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
# This sheet stores every value as text, even the numeric-looking ones
# (Magnitud/Latitud/Longitud/Profundidad such as "4.0", "7.0", "9.0").
# A plain `$cell.Value = "4.0"` on a General-formatted cell gets silently
# coerced to the number 4, which is NOT what we want here. To type a
# "numeric-looking" value while keeping it as genuine text (and without
# leaving a Text number-format on the destination cell), we stage the
# value in a scratch cell that IS formatted as Text, copy it, and paste
# only the *value* into the destination - that bypasses re-parsing the
# string as a number and leaves the destination cell's own formatting
# (General, no explicit style) untouched.
$scratch = $ws.Cells.Item(200, 200)
$scratch.NumberFormat = "@"

function Set-TextValue($row, $col, $val) {
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
}

# Plain text values (dates/times/labels) are never mistaken for numbers
# by this data, so a direct assignment already keeps them as text.
function Set-Text($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Row 5
Set-Text 5 1 "24/02/2020"
Set-Text 5 2 "02:02:08"
Set-TextValue 5 3 "4.0"
Set-Text 5 4 "CHOQUE_PLACAS"
Set-Text 5 5 "hola"
Set-TextValue 5 6 "7.0"
Set-TextValue 5 7 "9.0"
Set-TextValue 5 8 "4.0"
Set-Text 5 9 "LIMON, adios"

# Row 6
Set-Text 6 1 "24/02/2020"
Set-Text 6 2 "02:02:08"
Set-TextValue 6 3 "4.0"
Set-Text 6 4 "SUBDUCCION_PLACA"
Set-Text 6 5 "hola"
Set-TextValue 6 6 "7.0"
Set-TextValue 6 7 "9.0"
Set-TextValue 6 8 "4.0"
Set-Text 6 9 "HEREDIA, adios"

# Row 7
Set-Text 7 1 "24/02/2021"
Set-Text 7 2 "02:02:08"
Set-TextValue 7 3 "4.0"
Set-Text 7 4 "CHOQUE_PLACAS"
Set-Text 7 5 "hola"
Set-TextValue 7 6 "7.0"
Set-TextValue 7 7 "9.0"
Set-TextValue 7 8 "4.0"
Set-Text 7 9 "ALAJUELA, adios"

# Row 8
Set-Text 8 1 "24/02/2040"
Set-Text 8 2 "02:02:50"
Set-TextValue 8 3 "5.0"
Set-Text 8 4 "SUBDUCCION_PLACA"
Set-Text 8 5 "holo"
Set-TextValue 8 6 "8.0"
Set-TextValue 8 7 "8.0"
Set-TextValue 8 8 "7.0"
Set-Text 8 9 "PUNTARENAS, hola"

# Row 9
Set-Text 9 1 "24/02/2020"
Set-Text 9 2 "02:02:08"
Set-TextValue 9 3 "4.0"
Set-Text 9 4 "SUBDUCCION_PLACA"
Set-Text 9 5 "hola"
Set-TextValue 9 6 "7.0"
Set-TextValue 9 7 "9.0"
Set-TextValue 9 8 "4.0"
Set-Text 9 9 "SAN_JOSE, adios"

# Clean up the scratch cell so it doesn't leave stray data/dimension.
$scratch.Clear()
